$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efna5"
$ws.Cells.Item(2, 3).Value = "Ephb1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.172256
$ws.Cells.Item(2, 8).Value = 0.516768
$ws.Cells.Item(2, 9).Value = 0.1007998459820299
$ws.Cells.Item(2, 10).Value = 0.1116918428350044
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.875567333333333
$ws.Cells.Item(2, 14).Value = 5.626702
$ws.Cells.Item(2, 15).Value = 0.9147184316015459
$ws.Cells.Item(2, 16).Value = 0.9397527619538806
$ws.Cells.Item(2, 17).Value = 0.3230777265706666
$ws.Cells.Item(2, 18).Value = 2.907699539136
$ws.Cells.Item(2, 19).Value = 0.09220347702235974
$ws.Cells.Item(2, 20).Value = 0.1049627177919141

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efna5"
$ws.Cells.Item(3, 3).Value = "Ephb1"
$ws.Cells.Item(3, 4).Value = "M2"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.172256
$ws.Cells.Item(3, 8).Value = 0.516768
$ws.Cells.Item(3, 9).Value = 0.1007998459820299
$ws.Cells.Item(3, 10).Value = 0.1116918428350044
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.010998
$ws.Cells.Item(3, 14).Value = 0.032994
$ws.Cells.Item(3, 15).Value = 0.005363749481003509
$ws.Cells.Item(3, 16).Value = 0.005510546431623061
$ws.Cells.Item(3, 17).Value = 0.001894471488
$ws.Cells.Item(3, 18).Value = 0.017050243392
$ws.Cells.Item(3, 19).Value = 0.0005406651215713464
$ws.Cells.Item(3, 20).Value = 0.0006154830859758373

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efna5"
$ws.Cells.Item(4, 3).Value = "Ephb1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.172256
$ws.Cells.Item(4, 8).Value = 0.516768
$ws.Cells.Item(4, 9).Value = 0.1007998459820299
$ws.Cells.Item(4, 10).Value = 0.1116918428350044
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.163866
$ws.Cells.Item(4, 14).Value = 0.327732
$ws.Cells.Item(4, 15).Value = 0.07991781891745053
$ws.Cells.Item(4, 16).Value = 0.05473669161449624
$ws.Cells.Item(4, 17).Value = 0.028226901696
$ws.Cells.Item(4, 18).Value = 0.169361410176
$ws.Cells.Item(4, 19).Value = 0.008055703838098767
$ws.Cells.Item(4, 20).Value = 0.006113641957114418

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Efna5"
$ws.Cells.Item(5, 3).Value = "Ephb1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.036691
$ws.Cells.Item(5, 8).Value = 3.110073
$ws.Cells.Item(5, 9).Value = 0.6066453019398833
$ws.Cells.Item(5, 10).Value = 0.6721967782861762
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.875567333333333
$ws.Cells.Item(5, 14).Value = 5.626702
$ws.Cells.Item(5, 15).Value = 0.9147184316015459
$ws.Cells.Item(5, 16).Value = 0.9397527619538806
$ws.Cells.Item(5, 17).Value = 1.944383774360667
$ws.Cells.Item(5, 18).Value = 17.499453969246
$ws.Cells.Item(5, 19).Value = 0.5549096391288962
$ws.Cells.Item(5, 20).Value = 0.6316987789709344

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Efna5"
$ws.Cells.Item(6, 3).Value = "Ephb1"
$ws.Cells.Item(6, 4).Value = "M2"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.036691
$ws.Cells.Item(6, 8).Value = 3.110073
$ws.Cells.Item(6, 9).Value = 0.6066453019398833
$ws.Cells.Item(6, 10).Value = 0.6721967782861762
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.010998
$ws.Cells.Item(6, 14).Value = 0.032994
$ws.Cells.Item(6, 15).Value = 0.005363749481003509
$ws.Cells.Item(6, 16).Value = 0.005510546431623061
$ws.Cells.Item(6, 17).Value = 0.011401527618
$ws.Cells.Item(6, 18).Value = 0.102613748562
$ws.Cells.Item(6, 19).Value = 0.003253893423433266
$ws.Cells.Item(6, 20).Value = 0.003704171557933406

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efna5"
$ws.Cells.Item(7, 3).Value = "Ephb1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.036691
$ws.Cells.Item(7, 8).Value = 3.110073
$ws.Cells.Item(7, 9).Value = 0.6066453019398833
$ws.Cells.Item(7, 10).Value = 0.6721967782861762
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.163866
$ws.Cells.Item(7, 14).Value = 0.327732
$ws.Cells.Item(7, 15).Value = 0.07991781891745053
$ws.Cells.Item(7, 16).Value = 0.05473669161449624
$ws.Cells.Item(7, 17).Value = 0.169878407406
$ws.Cells.Item(7, 18).Value = 1.019270444436
$ws.Cells.Item(7, 19).Value = 0.0484817693875537
$ws.Cells.Item(7, 20).Value = 0.03679382775730833

$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Efna5"
$ws.Cells.Item(8, 3).Value = "Ephb1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.4999445
$ws.Cells.Item(8, 8).Value = 0.999889
$ws.Cells.Item(8, 9).Value = 0.2925548520780869
$ws.Cells.Item(8, 10).Value = 0.2161113788788194
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.875567333333333
$ws.Cells.Item(8, 14).Value = 5.626702
$ws.Cells.Item(8, 15).Value = 0.9147184316015459
$ws.Cells.Item(8, 16).Value = 0.9397527619538806
$ws.Cells.Item(8, 17).Value = 0.9376795726796667
$ws.Cells.Item(8, 18).Value = 5.626077436078
$ws.Cells.Item(8, 19).Value = 0.2676053154502899
$ws.Cells.Item(8, 20).Value = 0.203091265191032

$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Efna5"
$ws.Cells.Item(9, 3).Value = "Ephb1"
$ws.Cells.Item(9, 4).Value = "M2"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.4999445
$ws.Cells.Item(9, 8).Value = 0.999889
$ws.Cells.Item(9, 9).Value = 0.2925548520780869
$ws.Cells.Item(9, 10).Value = 0.2161113788788194
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.010998
$ws.Cells.Item(9, 14).Value = 0.032994
$ws.Cells.Item(9, 15).Value = 0.005363749481003509
$ws.Cells.Item(9, 16).Value = 0.005510546431623061
$ws.Cells.Item(9, 17).Value = 0.005498389611000001
$ws.Cells.Item(9, 18).Value = 0.03299033766600001
$ws.Cells.Item(9, 19).Value = 0.001569190935998897
$ws.Cells.Item(9, 20).Value = 0.001190891787713817

$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Efna5"
$ws.Cells.Item(10, 3).Value = "Ephb1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.4999445
$ws.Cells.Item(10, 8).Value = 0.999889
$ws.Cells.Item(10, 9).Value = 0.2925548520780869
$ws.Cells.Item(10, 10).Value = 0.2161113788788194
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.163866
$ws.Cells.Item(10, 14).Value = 0.327732
$ws.Cells.Item(10, 15).Value = 0.07991781891745053
$ws.Cells.Item(10, 16).Value = 0.05473669161449624
$ws.Cells.Item(10, 17).Value = 0.08192390543700001
$ws.Cells.Item(10, 18).Value = 0.327695621748
$ws.Cells.Item(10, 19).Value = 0.02338034569179808
$ws.Cells.Item(10, 20).Value = 0.01182922190007349
